# Auto-generated COM-interop script applying the Garuda_Profits.xlsx diff
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (ALC)
$ws.Range("H15").Value = 15926.029
$ws.Range("I15").Value = 15926.029
$ws.Range("K15").Value = 47778.087
$ws.Range("M15").Value = -47609.087

# Row 40 (ALC)
$ws.Range("H40").Value = 1552.3334
$ws.Range("I40").Value = 1114.1428
$ws.Range("J40").Value = 1771.4286
$ws.Range("K40").Value = 1114.1428
$ws.Range("L40").Value = 1771.4286
$ws.Range("M40").Value = -939.1428000000001
$ws.Range("N40").Value = -2121.4286

# Row 62 (ALC)
$ws.Range("H62").Value = 6201.6665
$ws.Range("I62").Value = 6802.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 6802.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -6178.5
$ws.Range("N62").Value = -6248

# Row 65 (ALC)
$ws.Range("H65").Value = 6201.6665
$ws.Range("I65").Value = 6802.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 34012.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -30892.5
$ws.Range("N65").Value = -31240

# Row 98 (ALC)
$ws.Range("H98").Value = 30829.363
$ws.Range("I98").Value = 464.28
$ws.Range("J98").Value = 125720.25
$ws.Range("K98").Value = 464.28
$ws.Range("L98").Value = 125720.25
$ws.Range("M98").Value = 1033.72
$ws.Range("N98").Value = -128716.25

# Row 116 (ALC)
$ws.Range("H116").Value = 1980.7693
$ws.Range("J116").Value = 1980.7693
$ws.Range("L116").Value = 1980.7693
$ws.Range("N116").Value = -8864.7693

# Row 122 (ALC)
$ws.Range("H122").Value = 30829.363
$ws.Range("I122").Value = 464.28
$ws.Range("J122").Value = 125720.25
$ws.Range("K122").Value = 1392.84
$ws.Range("L122").Value = 377160.75
$ws.Range("M122").Value = 1057.16
$ws.Range("N122").Value = -382060.75

# Row 125 (ALC)
$ws.Range("H125").Value = 100001040
$ws.Range("I125").Value = 200000660
$ws.Range("J125").Value = 1427.2
$ws.Range("K125").Value = 1800005940
$ws.Range("L125").Value = 12844.8
$ws.Range("M125").Value = -1800003480
$ws.Range("N125").Value = -17764.8

# Row 132 (ALC)
$ws.Range("H132").Value = 3323631.8
$ws.Range("I132").Value = 4082968.8
$ws.Range("K132").Value = 12248906.4
$ws.Range("M132").Value = -12246376.4

# Row 137 (ALC)
$ws.Range("H137").Value = 1331.75
$ws.Range("I137").Value = 986.9167
$ws.Range("J137").Value = 3400.75
$ws.Range("K137").Value = 2960.7501
$ws.Range("L137").Value = 10202.25
$ws.Range("M137").Value = -410.7501000000002
$ws.Range("N137").Value = -15302.25

# Row 141 (ALC)
$ws.Range("H141").Value = 2761.8076
$ws.Range("I141").Value = 1983
$ws.Range("J141").Value = 3379.4827
$ws.Range("K141").Value = 5949
$ws.Range("L141").Value = 10138.4481
$ws.Range("M141").Value = -769
$ws.Range("N141").Value = -20498.4481

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 62501132
$ws.Range("I2").Value = 111111740
$ws.Range("J2").Value = 1770.8572
$ws.Range("K2").Value = 111111740
$ws.Range("L2").Value = 1770.8572
$ws.Range("M2").Value = -111111627
$ws.Range("N2").Value = -1996.8572

# Row 45 (ARM)
$ws.Range("H45").Value = 15152852
$ws.Range("I45").Value = 33334386
$ws.Range("J45").Value = 1574.5834
$ws.Range("K45").Value = 33334386
$ws.Range("L45").Value = 1574.5834
$ws.Range("M45").Value = -33334009
$ws.Range("N45").Value = -2328.5834

# Row 110 (ARM)
$ws.Range("H110").Value = 2063.818
$ws.Range("I110").Value = 2189.111
$ws.Range("K110").Value = 2189.111
$ws.Range("M110").Value = -144.1109999999999

# Row 116 (ARM)
$ws.Range("H116").Value = 62501132
$ws.Range("I116").Value = 111111740
$ws.Range("J116").Value = 1770.8572
$ws.Range("K116").Value = 111111740
$ws.Range("L116").Value = 1770.8572
$ws.Range("M116").Value = -111109446
$ws.Range("N116").Value = -6358.8572

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 62501132
$ws.Range("I3").Value = 111111740
$ws.Range("J3").Value = 1770.8572
$ws.Range("K3").Value = 111111740
$ws.Range("L3").Value = 1770.8572
$ws.Range("M3").Value = -111111626
$ws.Range("N3").Value = -1998.8572

# Row 22 (BSM)
$ws.Range("H22").Value = 526.73334
$ws.Range("I22").Value = 526.73334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 526.73334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -353.73334
$ws.Range("N22").ClearContents()

# Row 99 (BSM)
$ws.Range("H99").Value = 2212.2
$ws.Range("I99").Value = 2050
$ws.Range("J99").Value = 2455.5
$ws.Range("K99").Value = 2050
$ws.Range("L99").Value = 2455.5
$ws.Range("M99").Value = -552
$ws.Range("N99").Value = -5451.5

# Row 105 (BSM)
$ws.Range("H105").Value = 2182.2354
$ws.Range("I105").Value = 1757
$ws.Range("J105").Value = 2479.9
$ws.Range("K105").Value = 1757
$ws.Range("L105").Value = 2479.9
$ws.Range("M105").Value = -10
$ws.Range("N105").Value = -5973.9

# Row 107 (BSM)
$ws.Range("H107").Value = 1986.2307
$ws.Range("I107").Value = 1881.4445
$ws.Range("J107").Value = 2222
$ws.Range("K107").Value = 1881.4445
$ws.Range("L107").Value = 2222
$ws.Range("M107").Value = 38.55549999999994
$ws.Range("N107").Value = -6062

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 861.6087
$ws.Range("I16").Value = 801.0625
$ws.Range("K16").Value = 801.0625
$ws.Range("M16").Value = -514.0625

# Row 31 (CRP)
$ws.Range("H31").Value = 2528127
$ws.Range("I31").Value = 2045.1714
$ws.Range("K31").Value = 2045.1714
$ws.Range("M31").Value = -1750.1714

# Row 34 (CRP)
$ws.Range("H34").Value = 2528127
$ws.Range("I34").Value = 2045.1714
$ws.Range("K34").Value = 2045.1714
$ws.Range("M34").Value = -1843.1714

# Row 110 (CRP)
$ws.Range("H110").Value = 29950.334
$ws.Range("J110").Value = 29950.334
$ws.Range("L110").Value = 29950.334
$ws.Range("N110").Value = -38130.334

# Row 112 (CRP)
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

# Row 113 (CRP)
$ws.Range("H113").Value = 861.6087
$ws.Range("I113").Value = 801.0625
$ws.Range("K113").Value = 801.0625
$ws.Range("M113").Value = 1368.9375

# Row 122 (CRP)
$ws.Range("H122").Value = 776.2941
$ws.Range("I122").Value = 767.9231
$ws.Range("J122").Value = 803.5
$ws.Range("K122").Value = 2303.7693
$ws.Range("L122").Value = 2410.5
$ws.Range("M122").Value = 146.2307000000001
$ws.Range("N122").Value = -7310.5

# Row 132 (CRP)
$ws.Range("H132").Value = 2365.862
$ws.Range("I132").Value = 2657
$ws.Range("J132").Value = 2094.1333
$ws.Range("K132").Value = 7971
$ws.Range("L132").Value = 6282.3999
$ws.Range("M132").Value = -5441
$ws.Range("N132").Value = -11342.3999

$ws = $wb.Worksheets.Item("CUL")
# Row 69 (CUL)
$ws.Range("H69").Value = 914.7368
$ws.Range("I69").Value = 460
$ws.Range("J69").Value = 1000
$ws.Range("K69").Value = 1380
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = -569
$ws.Range("N69").Value = -4622

# Row 72 (CUL)
$ws.Range("H72").Value = 914.7368
$ws.Range("I72").Value = 460
$ws.Range("J72").Value = 1000
$ws.Range("K72").Value = 4140
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = -84
$ws.Range("N72").Value = -17112

# Row 122 (CUL)
$ws.Range("H122").Value = 1077.2354
$ws.Range("I122").Value = 396.42856
$ws.Range("J122").Value = 1553.8
$ws.Range("K122").Value = 3567.85704
$ws.Range("L122").Value = 13984.2
$ws.Range("M122").Value = -1117.85704
$ws.Range("N122").Value = -18884.2

# Row 131 (CUL)
$ws.Range("H131").Value = 2060404.8
$ws.Range("J131").Value = 2526070
$ws.Range("L131").Value = 7578210
$ws.Range("N131").Value = -7588290

# Row 140 (CUL)
$ws.Range("H140").Value = 2305.4736
$ws.Range("I140").Value = 1912.5883
$ws.Range("J140").Value = 5645
$ws.Range("K140").Value = 5737.7649
$ws.Range("L140").Value = 16935
$ws.Range("M140").Value = -557.7649000000001
$ws.Range("N140").Value = -27295

$ws = $wb.Worksheets.Item("GSM")
# Row 52 (GSM)
$ws.Range("H52").Value = 11011
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 102 (GSM)
$ws.Range("H102").Value = 1158.7142
$ws.Range("I102").Value = 1098.9375
$ws.Range("J102").Value = 1350
$ws.Range("K102").Value = 1098.9375
$ws.Range("L102").Value = 1350
$ws.Range("M102").Value = 523.0625
$ws.Range("N102").Value = -4594

# Row 113 (GSM)
$ws.Range("H113").Value = 15625806
$ws.Range("I113").Value = 27778470
$ws.Range("J113").Value = 951.8570999999999
$ws.Range("K113").Value = 27778470
$ws.Range("L113").Value = 951.8570999999999
$ws.Range("M113").Value = -27776300
$ws.Range("N113").Value = -5291.8571

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 733
$ws.Range("I16").Value = 349.75
$ws.Range("J16").Value = 1499.5
$ws.Range("K16").Value = 349.75
$ws.Range("L16").Value = 1499.5
$ws.Range("M16").Value = -179.75
$ws.Range("N16").Value = -1839.5

# Row 22 (LTW)
$ws.Range("H22").Value = 628
$ws.Range("I22").Value = 718
$ws.Range("J22").Value = 520
$ws.Range("K22").Value = 718
$ws.Range("L22").Value = 520
$ws.Range("M22").Value = -423
$ws.Range("N22").Value = -1110

# Row 27 (LTW)
$ws.Range("H27").Value = 628
$ws.Range("I27").Value = 718
$ws.Range("J27").Value = 520
$ws.Range("K27").Value = 718
$ws.Range("L27").Value = 520
$ws.Range("M27").Value = -611
$ws.Range("N27").Value = -734

# Row 81 (LTW)
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 10000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -9002
$ws.Range("N81").ClearContents()

# Row 84 (LTW)
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 10000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -25008
$ws.Range("N84").ClearContents()

# Row 132 (LTW)
$ws.Range("H132").Value = 1682.6285
$ws.Range("I132").Value = 1303.9166
$ws.Range("J132").Value = 2508.9092
$ws.Range("K132").Value = 3911.7498
$ws.Range("L132").Value = 7526.7276
$ws.Range("M132").Value = -1381.7498
$ws.Range("N132").Value = -12586.7276

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 1064.5491
$ws.Range("I132").Value = 885.35895
$ws.Range("K132").Value = 2656.07685
$ws.Range("M132").Value = -126.0768500000004

# Row 136 (WVR)
$ws.Range("H136").Value = 3116.8333
$ws.Range("I136").Value = 3106.0588
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 9318.1764
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -6768.1764
$ws.Range("N136").Value = -15000
